$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "Kaustubh Kale, a software developer with 5+ years of experience in Java, Python, and C++, expresses interest in the advertised position. He highlights his skills, experience, and desire to contribute to the team, attaching his resume and cover letter for review."
$ws.Range("F2").Value = "Kaustubh Kale, un desarrollador de software con más de 5 años de experiencia en Java, Python y C++, expresa su interés en el puesto anunciado. Destaca sus habilidades, experiencia y deseo de contribuir al equipo, adjuntando su currículum vitae y carta de presentación para su revisión."

$ws.Range("E3").Value = "Tejas requests a meeting this week with Kaustubh to discuss the project timeline, ensuring alignment on milestones and deliverables before the deadline."
$ws.Range("F3").Value = "Tejas solicita una reunión esta semana con Kaustubh para discutir el cronograma del proyecto, asegurando la alineación en los hitos y entregables antes de la fecha límite."

$ws.Range("E4").Value = "The email requests branding guidelines, high-resolution images, and copy for the marketing campaign by Friday to begin the design phase and meet deadlines."
$ws.Range("F4").Value = "El correo electrónico solicita las directrices de marca, imágenes de alta resolución y el texto para la campaña de marketing para el viernes para comenzar la fase de diseño y cumplir con las fechas límite."

$ws.Range("E5").Value = "Saloni Goyal inquires about international student admission requirements, including documentation, English proficiency tests, deadlines, and scholarship opportunities."
$ws.Range("F5").Value = "Saloni Goyal consulta sobre los requisitos de admisión para estudiantes internacionales, incluyendo la documentación, las pruebas de dominio del inglés, las fechas límite y las oportunidades de becas."

$ws.Range("E6").Value = "Customer Niraj Wadhe reports a defective blender (Order #12345) and requests instructions for return and replacement, emphasizing the urgent need for a functional blender."
$ws.Range("F6").Value = "El cliente Niraj Wadhe reporta una licuadora defectuosa (Orden #12345) y solicita instrucciones para la devolución y reemplazo, enfatizando la necesidad urgente de una licuadora funcional."

$ws.Range("E7").Value = "Tejas expresses gratitude to Kaustubh for his exceptional contributions to the project, highlighting his expertise, dedication, and willingness to go the extra mile."
$ws.Range("F7").Value = "Tejas expresa su agradecimiento a Kaustubh por sus excepcionales contribuciones al proyecto, destacando su experiencia, dedicación y disposición a esforzarse más allá de lo esperado."

$ws.Range("E8").Value = "You are invited to the Annual Gala on July 20th at the Grand Ballroom, featuring dinner, entertainment, and an auction to support charitable initiatives. RSVP by July 10th."
$ws.Range("F8").Value = "Estás invitado a la Gala Anual el 20 de julio en el Gran Salón de Baile, que incluirá cena, entretenimiento y una subasta para apoyar iniciativas benéficas. Confirma tu asistencia antes del 10 de julio."

$ws.Range("E9").Value = "Niraj Wadhe admires Tejas' article on industry trends and seeks collaboration, proposing a call next week to explore potential partnerships leveraging Tejas' expertise."
$ws.Range("F9").Value = "Niraj Wadhe admira el artículo de Tejas sobre las tendencias de la industria y busca colaborar, proponiendo una llamada la próxima semana para explorar posibles asociaciones aprovechando la experiencia de Tejas."

$ws.Range("E10").Value = "This email reminds you of your appointment with Dr. Smith on June 30th at 10:00 AM. Please arrive 10 minutes early and bring medical records. Contact the office if you need to reschedule or have questions."
$ws.Range("F10").Value = "Este correo electrónico es para recordarte tu cita con el Dr. Smith el 30 de junio a las 10:00 AM. Por favor, llega 10 minutos antes y trae tus registros médicos. Contacta a la oficina si necesitas reprogramar la cita o tienes alguna pregunta."

$ws.Range("E11").Value = "This June newsletter highlights new products, special discounts, and expert articles. Click the link to explore exclusive offers and stay updated on future innovations."
$ws.Range("F11").Value = "Este boletín de junio destaca nuevos productos, descuentos especiales y artículos de expertos. Haga clic en el enlace para explorar ofertas exclusivas y mantenerse actualizado sobre futuras innovaciones."
